{"js": "// Apply the edits described by the diff: update application numbers/dates,\n// company name, service location, signatory names, authority basis,\n// survey/certificate descriptions, and monetary amounts.\n//\n// Each pair is applied as a global (all-occurrences) replace, which matches\n// the diff exactly (strings appearing twice in the template are replaced\n// identically both times; unique strings appear once).\nconst replacements = [\n  [\"2409865\", \"3423567\"],\n  [\"05.04.2024\", \"01.08.2024\"],\n  [\n    \"\u0410\u043a\u0446\u0438\u043e\u043d\u0435\u0440\u043d\u043e\u0435 \u043e\u0431\u0449\u0435\u0441\u0442\u0432\u043e \u00ab\u0410\u0440\u043a\u0442\u0438\u0447\u0435\u0441\u043a\u0438\u0435 \u043c\u043e\u0440\u0441\u043a\u0438\u0435 \u0438\u043d\u0436\u0435\u043d\u0435\u0440\u043d\u043e-\u0433\u0435\u043e\u043b\u043e\u0433\u0438\u0447\u0435\u0441\u043a\u0438\u0435 \u044d\u043a\u0441\u043f\u0435\u0434\u0438\u0446\u0438\u0438\u00bb\",\n    \"\u041e\u0431\u0449\u0435\u0441\u0442\u0432\u043e \u0441 \u043e\u0433\u0440\u0430\u043d\u0438\u0447\u0435\u043d\u043d\u043e\u0439 \u043e\u0442\u0432\u0435\u0442\u0441\u0442\u0432\u0435\u043d\u043d\u043e\u0441\u0442\u044c\u044e \u00ab\u0421\u041f\u0415\u0426 \u0411\u0410\u041b\u0422\u00bb\",\n  ],\n  [\"\u0421\u0432\u0435\u0442\u043b\u044b\u0439 \", \"\u0421\u0432\u0435\u0442\u043b\u044b\u0439, \u041a\u0430\u043b\u0438\u043d\u0438\u043d\u0433\u0440\u0430\u0434\u0441\u043a\u0430\u044f \u043e\u0431\u043b\u0430\u0441\u0442\u044c, \u0420\u043e\u0441\u0441\u0438\u044f\"],\n  [\"\u0421\u0443\u043f\u0435\u0440\u0438\u043d\u0442\u0435\u043d\u0434\u0430\u043d\u0442 \u041c\u0443\u0445\u0438\u043d \u041a. \u0410.\", \"\u0433\u0435\u043d\u0435\u0440\u0430\u043b\u044c\u043d\u044b\u0439 \u0434\u0438\u0440\u0435\u043a\u0442\u043e\u0440 \u041c\u0430\u0440\u043a\u043e\u0432\u0430 \u041d. \u0412.\"],\n  [\"\u0414\u043e\u0432\u0435\u0440\u0435\u043d\u043d\u043e\u0441\u0442\u0438 \u2116  \u043e\u0442 -- \", \"\u0423\u0441\u0442\u0430\u0432\u0430\"],\n  [\n    \"\u0415\u0436\u0435\u0433\u043e\u0434\u043d\u043e\u0435 \u043e\u0441\u0432\u0438\u0434\u0435\u0442\u0435\u043b\u044c\u0441\u0442\u0432\u043e\u0432\u0430\u043d\u0438\u0435 \u0421\u0421\u041f \u2116 24.42.03.00765.121 \u043e\u0442 02.05.2021 \",\n    \"\u0415\u0436\u0435\u0433\u043e\u0434\u043d\u043e\u0435 \u043e\u0441\u0432\u0438\u0434\u0435\u0442\u0435\u043b\u044c\u0441\u0442\u0432\u043e\u0432\u0430\u043d\u0438\u0435 \u043e\u0444\u0444\u0448\u043e\u0440\u043d\u044b\u0445 \u043a\u043e\u043d\u0442\u0435\u0439\u043d\u0435\u0440\u043e\u0432 \u043d\u0430 \u0441\u043e\u043e\u0442\u0432\u0435\u0442\u0441\u0442\u0432\u0438\u0435 \u0442\u0440\u0435\u0431\u043e\u0432\u0430\u043d\u0438\u044f\u043c \u041a\u0411\u041a - 1 \u0448\u0442.\",\n  ],\n  [\n    \"\u0421\u0432\u0438\u0434\u0435\u0442\u0435\u043b\u044c\u0441\u0442\u0432\u043e \u0444. 7.1.27 \u2116 24.02.42.00987.121 \u043e\u0442 --\",\n    \"\u0421\u0432\u0438\u0434\u0435\u0442\u0435\u043b\u044c\u0441\u0442\u0432\u043e \u0444. 7.1.27 \u2116 2422442432342424 \u043e\u0442 --\",\n  ],\n  [\n    \"15 600,00 p. (\u043f\u044f\u0442\u043d\u0430\u0434\u0446\u0430\u0442\u044c \u0442\u044b\u0441\u044f\u0447 \u0448\u0435\u0441\u0442\u044c\u0441\u043e\u0442 \u0440\u0443\u0431\u043b\u0435\u0439 00 \u043a\u043e\u043f\u0435\u0435\u043a)\",\n    \"100 000,00 p. (\u0441\u0442\u043e \u0442\u044b\u0441\u044f\u0447 \u0440\u0443\u0431\u043b\u0435\u0439 00 \u043a\u043e\u043f\u0435\u0435\u043a)\",\n  ],\n  [\n    \"3 120,00 p. (\u0442\u0440\u0438 \u0442\u044b\u0441\u044f\u0447\u0438 \u0441\u0442\u043e \u0434\u0432\u0430\u0434\u0446\u0430\u0442\u044c \u0440\u0443\u0431\u043b\u0435\u0439 00 \u043a\u043e\u043f\u0435\u0435\u043a)\",\n    \"20 000,00 p. (\u0434\u0432\u0430\u0434\u0446\u0430\u0442\u044c \u0442\u044b\u0441\u044f\u0447 \u0440\u0443\u0431\u043b\u0435\u0439 00 \u043a\u043e\u043f\u0435\u0435\u043a)\",\n  ],\n  [\n    \"18 720,00 p. (\u0432\u043e\u0441\u0435\u043c\u043d\u0430\u0434\u0446\u0430\u0442\u044c \u0442\u044b\u0441\u044f\u0447 \u0441\u0435\u043c\u044c\u0441\u043e\u0442 \u0434\u0432\u0430\u0434\u0446\u0430\u0442\u044c \u0440\u0443\u0431\u043b\u0435\u0439 00 \u043a\u043e\u043f\u0435\u0435\u043a)\",\n    \"120 000,00 p. (\u0441\u0442\u043e \u0434\u0432\u0430\u0434\u0446\u0430\u0442\u044c \u0442\u044b\u0441\u044f\u0447 \u0440\u0443\u0431\u043b\u0435\u0439 00 \u043a\u043e\u043f\u0435\u0435\u043a)\",\n  ],\n  [\"\u041a. \u0410. \u041c\u0443\u0445\u0438\u043d\", \"\u041d. \u0412. \u041c\u0430\u0440\u043a\u043e\u0432\u0430\"],\n];\n\nconst body = context.document.body;\n\nfor (const [searchText, newText] of replacements) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the edits described by the diff: update application numbers/dates,\n# company name, service location, signatory names, authority basis,\n# survey/certificate descriptions, and monetary amounts.\n#\n# Each pair is applied via Find/Replace across the whole document body\n# (wdReplaceAll), which matches the diff exactly (strings appearing twice\n# in the template are replaced identically both times; unique strings\n# appear once).\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2409865\", \"3423567\"),\n    @(\"05.04.2024\", \"01.08.2024\"),\n    @(\"\u0410\u043a\u0446\u0438\u043e\u043d\u0435\u0440\u043d\u043e\u0435 \u043e\u0431\u0449\u0435\u0441\u0442\u0432\u043e \u00ab\u0410\u0440\u043a\u0442\u0438\u0447\u0435\u0441\u043a\u0438\u0435 \u043c\u043e\u0440\u0441\u043a\u0438\u0435 \u0438\u043d\u0436\u0435\u043d\u0435\u0440\u043d\u043e-\u0433\u0435\u043e\u043b\u043e\u0433\u0438\u0447\u0435\u0441\u043a\u0438\u0435 \u044d\u043a\u0441\u043f\u0435\u0434\u0438\u0446\u0438\u0438\u00bb\", \"\u041e\u0431\u0449\u0435\u0441\u0442\u0432\u043e \u0441 \u043e\u0433\u0440\u0430\u043d\u0438\u0447\u0435\u043d\u043d\u043e\u0439 \u043e\u0442\u0432\u0435\u0442\u0441\u0442\u0432\u0435\u043d\u043d\u043e\u0441\u0442\u044c\u044e \u00ab\u0421\u041f\u0415\u0426 \u0411\u0410\u041b\u0422\u00bb\"),\n    @(\"\u0421\u0432\u0435\u0442\u043b\u044b\u0439 \", \"\u0421\u0432\u0435\u0442\u043b\u044b\u0439, \u041a\u0430\u043b\u0438\u043d\u0438\u043d\u0433\u0440\u0430\u0434\u0441\u043a\u0430\u044f \u043e\u0431\u043b\u0430\u0441\u0442\u044c, \u0420\u043e\u0441\u0441\u0438\u044f\"),\n    @(\"\u0421\u0443\u043f\u0435\u0440\u0438\u043d\u0442\u0435\u043d\u0434\u0430\u043d\u0442 \u041c\u0443\u0445\u0438\u043d \u041a. \u0410.\", \"\u0433\u0435\u043d\u0435\u0440\u0430\u043b\u044c\u043d\u044b\u0439 \u0434\u0438\u0440\u0435\u043a\u0442\u043e\u0440 \u041c\u0430\u0440\u043a\u043e\u0432\u0430 \u041d. \u0412.\"),\n    @(\"\u0414\u043e\u0432\u0435\u0440\u0435\u043d\u043d\u043e\u0441\u0442\u0438 \u2116  \u043e\u0442 -- \", \"\u0423\u0441\u0442\u0430\u0432\u0430\"),\n    @(\"\u0415\u0436\u0435\u0433\u043e\u0434\u043d\u043e\u0435 \u043e\u0441\u0432\u0438\u0434\u0435\u0442\u0435\u043b\u044c\u0441\u0442\u0432\u043e\u0432\u0430\u043d\u0438\u0435 \u0421\u0421\u041f \u2116 24.42.03.00765.121 \u043e\u0442 02.05.2021 \", \"\u0415\u0436\u0435\u0433\u043e\u0434\u043d\u043e\u0435 \u043e\u0441\u0432\u0438\u0434\u0435\u0442\u0435\u043b\u044c\u0441\u0442\u0432\u043e\u0432\u0430\u043d\u0438\u0435 \u043e\u0444\u0444\u0448\u043e\u0440\u043d\u044b\u0445 \u043a\u043e\u043d\u0442\u0435\u0439\u043d\u0435\u0440\u043e\u0432 \u043d\u0430 \u0441\u043e\u043e\u0442\u0432\u0435\u0442\u0441\u0442\u0432\u0438\u0435 \u0442\u0440\u0435\u0431\u043e\u0432\u0430\u043d\u0438\u044f\u043c \u041a\u0411\u041a - 1 \u0448\u0442.\"),\n    @(\"\u0421\u0432\u0438\u0434\u0435\u0442\u0435\u043b\u044c\u0441\u0442\u0432\u043e \u0444. 7.1.27 \u2116 24.02.42.00987.121 \u043e\u0442 --\", \"\u0421\u0432\u0438\u0434\u0435\u0442\u0435\u043b\u044c\u0441\u0442\u0432\u043e \u0444. 7.1.27 \u2116 2422442432342424 \u043e\u0442 --\"),\n    @(\"15 600,00 p. (\u043f\u044f\u0442\u043d\u0430\u0434\u0446\u0430\u0442\u044c \u0442\u044b\u0441\u044f\u0447 \u0448\u0435\u0441\u0442\u044c\u0441\u043e\u0442 \u0440\u0443\u0431\u043b\u0435\u0439 00 \u043a\u043e\u043f\u0435\u0435\u043a)\", \"100 000,00 p. (\u0441\u0442\u043e \u0442\u044b\u0441\u044f\u0447 \u0440\u0443\u0431\u043b\u0435\u0439 00 \u043a\u043e\u043f\u0435\u0435\u043a)\"),\n    @(\"3 120,00 p. (\u0442\u0440\u0438 \u0442\u044b\u0441\u044f\u0447\u0438 \u0441\u0442\u043e \u0434\u0432\u0430\u0434\u0446\u0430\u0442\u044c \u0440\u0443\u0431\u043b\u0435\u0439 00 \u043a\u043e\u043f\u0435\u0435\u043a)\", \"20 000,00 p. (\u0434\u0432\u0430\u0434\u0446\u0430\u0442\u044c \u0442\u044b\u0441\u044f\u0447 \u0440\u0443\u0431\u043b\u0435\u0439 00 \u043a\u043e\u043f\u0435\u0435\u043a)\"),\n    @(\"18 720,00 p. (\u0432\u043e\u0441\u0435\u043c\u043d\u0430\u0434\u0446\u0430\u0442\u044c \u0442\u044b\u0441\u044f\u0447 \u0441\u0435\u043c\u044c\u0441\u043e\u0442 \u0434\u0432\u0430\u0434\u0446\u0430\u0442\u044c \u0440\u0443\u0431\u043b\u0435\u0439 00 \u043a\u043e\u043f\u0435\u0435\u043a)\", \"120 000,00 p. (\u0441\u0442\u043e \u0434\u0432\u0430\u0434\u0446\u0430\u0442\u044c \u0442\u044b\u0441\u044f\u0447 \u0440\u0443\u0431\u043b\u0435\u0439 00 \u043a\u043e\u043f\u0435\u0435\u043a)\"),\n    @(\"\u041a. \u0410. \u041c\u0443\u0445\u0438\u043d\", \"\u041d. \u0412. \u041c\u0430\u0440\u043a\u043e\u0432\u0430\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
